$wb = $excel.ActiveWorkbook

# --- Sheet "Jaccard" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Jaccard")

# Insert two new rows before current row 4 ("Total"), shifting Residual/Total down
$ws1.Rows.Item(4).Insert()
$ws1.Rows.Item(4).Insert()

# Row 2: Model -> level, with new values
$ws1.Range("A2").Value = "level"
$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 0.7852985403925928
$ws1.Range("D2").Value = 0.06904394947671905
$ws1.Range("E2").Value = 1.283723165649595
$ws1.Range("F2").Value = 0.0062

# Row 3: Residual -> site, with new values
$ws1.Range("A3").Value = "site"
$ws1.Range("B3").Value = 3
$ws1.Range("C3").Value = 2.188545411289753
$ws1.Range("D3").Value = 0.1924183110400957
$ws1.Range("E3").Value = 2.385068691748097
$ws1.Range("F3").Value = 0.0001

# Row 4: new "site:plant" row
$ws1.Range("A4").Value = "site:plant"
$ws1.Range("B4").Value = 7
$ws1.Range("C4").Value = 3.20030123446554
$ws1.Range("D4").Value = 0.2813725295252116
$ws1.Range("E4").Value = 1.49471847617761
$ws1.Range("F4").Value = 0.0001

# Row 5: Residual (re-created, values unchanged from original row3)
$ws1.Range("A5").Value = "Residual"
$ws1.Range("B5").Value = 17
$ws1.Range("C5").Value = 5.199748490913387
$ws1.Range("D5").Value = 0.4571652099579737

# Row 6: Total (values unchanged, shifted from row4)
$ws1.Range("A6").Value = "Total"
$ws1.Range("B6").Value = 29
$ws1.Range("C6").Value = 11.37389367706127
$ws1.Range("D6").Value = 1

# --- Sheet "turnover" (sheet2) ---
$ws2 = $wb.Worksheets.Item("turnover")

$ws2.Rows.Item(4).Insert()
$ws2.Rows.Item(4).Insert()

$ws2.Range("A2").Value = "level"
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 0.4728651059065685
$ws2.Range("D2").Value = 0.04526884405202474
$ws2.Range("E2").Value = 0.9054736837328426
$ws2.Range("F2").Value = 0.7792

$ws2.Range("A3").Value = "site"
$ws2.Range("B3").Value = 3
$ws2.Range("C3").Value = 2.466256733236879
$ws2.Range("D3").Value = 0.2361024107184082
$ws2.Range("E3").Value = 3.148368729800303
$ws2.Range("F3").Value = 0.0001

$ws2.Range("A4").Value = "site:plant"
$ws2.Range("B4").Value = 7
$ws2.Range("C4").Value = 3.067634611377764
$ws2.Range("D4").Value = 0.2936741812759006
$ws2.Range("E4").Value = 1.678317604875803
$ws2.Range("F4").Value = 0.0001

$ws2.Range("A5").Value = "Residual"
$ws2.Range("B5").Value = 17
$ws2.Range("C5").Value = 4.438951095338217
$ws2.Range("D5").Value = 0.4249545639536665

$ws2.Range("A6").Value = "Total"
$ws2.Range("B6").Value = 29
$ws2.Range("C6").Value = 10.44570754585943
$ws2.Range("D6").Value = 1
